$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Swap the header text shown for columns B and D (KOR_GBR_841810 <-> CHN_GBR_841810)
$ws.Range("B1").Value = "CHN_GBR_841810"
$ws.Range("D1").Value = "KOR_GBR_841810"

# 2. Rows 12-47: only column D has a value; move it into column B and clear D.
for ($r = 12; $r -le 47; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    if ($dVal -ne "") {
        $ws.Cells.Item($r, 2).Value = $dVal
        $ws.Cells.Item($r, 4).ClearContents()
    }
}

# 3. Rows 48-167: both B and D have values; swap them.
for ($r = 48; $r -le 167; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 2).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $bVal
}

# 4. Rows 168-173: only column B has a value; move it into column D and clear B.
for ($r = 168; $r -le 173; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    if ($bVal -ne "") {
        $ws.Cells.Item($r, 4).Value = $bVal
        $ws.Cells.Item($r, 2).ClearContents()
    }
}
